# Update CSP benchmark results sheet: add cplex, gurobi, ortools results
# (shifts WANG1-3 rows earlier, CHL2 row to the end, and updates Runtime values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BENG01
$ws.Range("D2").Value = 0.213964613

# Row 3: BENG02
$ws.Range("D3").Value = 16.44030466

# Row 4: BENG03
$ws.Range("D4").Value = 11.2207408

# Row 5: BENG04
$ws.Range("D5").Value = 37.39759781

# Row 6: BENG05
$ws.Range("D6").Value = 70.94194152999999

# Row 7: BENG06
$ws.Range("D7").Value = 0.774997429

# Row 8: BENG07
$ws.Range("D8").Value = 13.67573711

# Row 9: BENG08
$ws.Range("D9").Value = 27.53228301

# Row 10: BENG09
$ws.Range("D10").Value = 72.41773309

# Row 11: BENG10
$ws.Range("D11").Value = 190.7715986

# Row 12: WANG1
$ws.Range("A12").Value = "WANG1"
$ws.Range("B12").Value = 6340
$ws.Range("C12").Value = 1193326
$ws.Range("D12").Value = 33.34636053
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = "TIMEOUT"

# Row 13: WANG2
$ws.Range("A13").Value = "WANG2"
$ws.Range("B13").Value = 6550
$ws.Range("C13").Value = 1115244
$ws.Range("D13").Value = 31.96047847
$ws.Range("E13").Value = 11
$ws.Range("F13").Value = "TIMEOUT"

# Row 14: WANG3
$ws.Range("A14").Value = "WANG3"
$ws.Range("B14").Value = 6592
$ws.Range("C14").Value = 1132506
$ws.Range("D14").Value = 11.4881725
$ws.Range("E14").Value = 11
$ws.Range("F14").Value = "TIMEOUT"

# Row 15: ngcut1
$ws.Range("A15").Value = "ngcut1"
$ws.Range("B15").Value = 319
$ws.Range("C15").Value = 2391
$ws.Range("D15").Value = 0.017789558
$ws.Range("E15").Value = 3

# Row 16: ngcut2
$ws.Range("A16").Value = "ngcut2"
$ws.Range("B16").Value = 805
$ws.Range("C16").Value = 10834
$ws.Range("D16").Value = 0.137762133
$ws.Range("E16").Value = 4

# Row 17: ngcut3
$ws.Range("A17").Value = "ngcut3"
$ws.Range("B17").Value = 1191
$ws.Range("C17").Value = 18366
$ws.Range("D17").Value = 0.150790223

# Row 18: ngcut4
$ws.Range("A18").Value = "ngcut4"
$ws.Range("B18").Value = 194
$ws.Range("C18").Value = 1318
$ws.Range("D18").Value = 0.004480143
$ws.Range("E18").Value = 2

# Row 19: ngcut5
$ws.Range("A19").Value = "ngcut5"
$ws.Range("B19").Value = 597
$ws.Range("C19").Value = 8025
$ws.Range("D19").Value = 0.024679722
$ws.Range("E19").Value = 3

# Row 20: ngcut6
$ws.Range("A20").Value = "ngcut6"
$ws.Range("B20").Value = 701
$ws.Range("C20").Value = 7549
$ws.Range("D20").Value = 0.122805094
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = "COMPLETE"

# Row 21: ngcut7
$ws.Range("A21").Value = "ngcut7"
$ws.Range("B21").Value = 356
$ws.Range("C21").Value = 1924
$ws.Range("D21").Value = 0.006977101
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "COMPLETE"

# Row 22: ngcut8
$ws.Range("A22").Value = "ngcut8"
$ws.Range("B22").Value = 676
$ws.Range("C22").Value = 8658
$ws.Range("D22").Value = 0.026513081
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = "COMPLETE"

# Row 23: ngcut9
$ws.Range("D23").Value = 0.180434531

# Row 24: ngcut10
$ws.Range("D24").Value = 0.156260539

# Row 25: ngcut11
$ws.Range("D25").Value = 0.140559098

# Row 26: ngcut12
$ws.Range("D26").Value = 0.385654115

# Row 27: cgcut1
$ws.Range("D27").Value = 0.031486081

# Row 28: cgcut2
$ws.Range("D28").Value = 0.318296757

# Row 29: cgcut3
$ws.Range("D29").Value = 9.334168258

# Row 30: A1
$ws.Range("D30").Value = 8.694348646

# Row 31: A2
$ws.Range("D31").Value = 32.14191192

# Row 32: A3
$ws.Range("D32").Value = 8.296645525000001

# Row 33: A4
$ws.Range("D33").Value = 3.330431542

# Row 34: A5
$ws.Range("D34").Value = 55.65743394

# Row 35: HH
$ws.Range("D35").Value = 0.424607828

# Row 36: CHL1
$ws.Range("D36").Value = 25.69063781

# Row 37: CHL3
$ws.Range("A37").Value = "CHL3"
$ws.Range("B37").Value = 11121
$ws.Range("C37").Value = 305340
$ws.Range("D37").Value = 2.293431980004243
$ws.Range("E37").Value = 1

# Row 38: CHL4
$ws.Range("A38").Value = "CHL4"
$ws.Range("B38").Value = 12062
$ws.Range("C38").Value = 287415
$ws.Range("D38").Value = 2.424728341000446

# Row 39: CHL5
$ws.Range("A39").Value = "CHL5"
$ws.Range("B39").Value = 1131
$ws.Range("C39").Value = 24783
$ws.Range("D39").Value = 0.1791259240053478
$ws.Range("E39").Value = 3

# Row 40: CHL6
$ws.Range("A40").Value = "CHL6"
$ws.Range("B40").Value = 21037
$ws.Range("C40").Value = 3989467
$ws.Range("D40").Value = 26.52379394799937
$ws.Range("E40").Value = 6
$ws.Range("F40").Value = "TIMEOUT"

# Row 41: CHL7
$ws.Range("A41").Value = "CHL7"
$ws.Range("B41").Value = 26109
$ws.Range("C41").Value = 6497355
$ws.Range("D41").Value = 71.77097530799801
$ws.Range("F41").Value = "COMPLETE"

# Row 42: Hchl1
$ws.Range("A42").Value = "Hchl1"
$ws.Range("B42").Value = 21037
$ws.Range("C42").Value = 3989467
$ws.Range("D42").Value = 26.48783655300213
$ws.Range("F42").Value = "TIMEOUT"

# Row 43: Hchl2
$ws.Range("A43").Value = "Hchl2"
$ws.Range("B43").Value = 26049
$ws.Range("C43").Value = 6470655
$ws.Range("D43").Value = 83.65945453500171
$ws.Range("F43").Value = "COMPLETE"

# Row 44: Hchl3s
$ws.Range("A44").Value = "Hchl3s"
$ws.Range("B44").Value = 14170
$ws.Range("C44").Value = 1350700
$ws.Range("D44").Value = 13.45013016500161
$ws.Range("E44").Value = 3

# Row 45: Hchl4s
$ws.Range("A45").Value = "Hchl4s"
$ws.Range("B45").Value = 7598
$ws.Range("C45").Value = 350770
$ws.Range("D45").Value = 3.603140851999342
$ws.Range("E45").Value = 2

# Row 46: Hchl5s
$ws.Range("A46").Value = "Hchl5s"
$ws.Range("B46").Value = 27280
$ws.Range("C46").Value = 4738380
$ws.Range("D46").Value = 47.51820178899652
$ws.Range("E46").Value = 4

# Row 47: Hchl6s
$ws.Range("A47").Value = "Hchl6s"
$ws.Range("B47").Value = 29417
$ws.Range("C47").Value = 6532942
$ws.Range("D47").Value = 65.66112957699806
$ws.Range("E47").Value = 5

# Row 48: Hchl7s
$ws.Range("A48").Value = "Hchl7s"
$ws.Range("B48").Value = 50118
$ws.Range("C48").Value = 20915127
$ws.Range("D48").Value = 297.2384830459996
$ws.Range("E48").Value = 7

# Row 49: Hchl8s
$ws.Range("A49").Value = "Hchl8s"
$ws.Range("B49").Value = 1617
$ws.Range("C49").Value = 17937
$ws.Range("D49").Value = 390.8667347299997
$ws.Range("E49").Value = 2
$ws.Range("F49").Value = "COMPLETE"

# Row 50: Hchl9
$ws.Range("A50").Value = "Hchl9"
$ws.Range("B50").Value = 19375
$ws.Range("C50").Value = 5451829
$ws.Range("D50").Value = 123.0014756040036
$ws.Range("E50").Value = 10

# Row 51: CHL2
$ws.Range("A51").Value = "CHL2"
$ws.Range("B51").Value = 2242
$ws.Range("C51").Value = 83220
$ws.Range("D51").Value = 0.2611286970000037
$ws.Range("E51").Value = 3

